$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the confidential disclosure date text in cell A7
$ws.Range("A7").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

# Update the numeric values in D2:E4
$ws.Range("D2").Value = 0.8449746302493792
$ws.Range("E2").Value = 0.01022102976874928

$ws.Range("D3").Value = 0.1550253697506208
$ws.Range("E3").Value = 0.01448467966573808

$ws.Range("E4").Value = 0.01088200367051706
